$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 63, shifting existing rows 63:74 down to 64:75.
$ws.Rows.Item(63).Insert()

# Populate the newly inserted row 63 with the new record.
$ws.Cells.Item(63, 1).Value = 5
$ws.Cells.Item(63, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(63, 3).Value = "Maule"
$ws.Cells.Item(63, 4).Value = 44876
$ws.Cells.Item(63, 5).Value = 7
$ws.Cells.Item(63, 6).Value = "Fruta"
$ws.Cells.Item(63, 7).Value = 100101
$ws.Cells.Item(63, 8).Value = "Berries"
$ws.Cells.Item(63, 9).Value = 100101001
$ws.Cells.Item(63, 10).Value = "Arándano (blue)"
$ws.Cells.Item(63, 11).Value = "Sin especificar"
$ws.Cells.Item(63, 12).Value = "Primera"
$ws.Cells.Item(63, 13).Value = 200
$ws.Cells.Item(63, 14).Value = 6000
$ws.Cells.Item(63, 15).Value = 6000
$ws.Cells.Item(63, 16).Value = 6000
$ws.Cells.Item(63, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(63, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(63, 19).Value = 3000
$ws.Cells.Item(63, 20).Value = 2
